$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet (sheet index 1): insert the new "2022-Q3" row at the top of
#    the data (row 2), push every existing row down by one, and append the
#    "2020-Q4" row that falls off the end as the new last row (row 9).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make column-A's index style (border/alignment) available on the new row 9
# by copying the formatting from row 8 (a single cheap cell copy).
$summary.Range("A8").Copy($summary.Range("A9"))

$summaryRows = @(
    @(0, "2022-Q3", 18, 3.49),
    @(1, "2022-Q2", 8, 1.1),
    @(2, "2022-Q1", 9, 0.9399999999999999),
    @(3, "2021-Q4", 5, 0.55),
    @(4, "2021-Q3", 3, 0.32),
    @(5, "2021-Q2", 1, 0),
    @(6, "2021-Q1", 2, 0.01),
    @(7, "2020-Q4", 2, 0.01)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right after "总计", pushing every
#    other quarter sheet down one slot (matches the target tab order).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Borrow the header/row formatting (fonts, borders, alignment) from the
# "2022-Q2" sheet so the new sheet matches the house style used by every
# other quarter tab.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("A1:H1").Copy($q3.Range("A1:H1"))
$template.Range("A2:H2").Copy($q3.Range("A2:H19"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# index, code, name, scale, stockPosition, positionRatio, marketValue, rank
$fundRows = @(
    @(0,  "009100", "安信稳健增利混合A",                       "113.23", "23.20", "0.96", "1.0870", 8),
    @(1,  "009101", "安信稳健增利混合C",                       "84.02",  "23.20", "0.96", "0.8066", 8),
    @(2,  "012609", "安信稳健汇利一年持有混合A",               "35.76",  "23.15", "0.87", "0.3111", 10),
    @(3,  "012256", "安信丰穗一年持有混合A",                   "24.55",  "24.64", "1.08", "0.2651", 9),
    @(4,  "008809", "安信民稳增长混合A",                       "9.52",   "49.26", "2.13", "0.2028", 9),
    @(5,  "009849", "安信稳健聚申一年持有期混合A",             "12.83",  "39.19", "1.48", "0.1899", 8),
    @(6,  "008810", "安信民稳增长混合C",                       "7.95",   "49.26", "2.13", "0.1693", 9),
    @(7,  "012610", "安信稳健汇利一年持有混合C",               "18.98",  "23.15", "0.87", "0.1651", 10),
    @(8,  "012250", "安信平衡增利混合A",                       "2.63",   "60.12", "2.21", "0.0581", 10),
    @(9,  "011651", "招商港股通核心精选股票A",                 "1.66",   "84.28", "2.80", "0.0465", 7),
    @(10, "012251", "安信平衡增利混合C",                       "2.10",   "60.12", "2.21", "0.0464", 10),
    @(11, "004266", "招商沪港深科技创新主题精选灵活配置混合A", "0.92",   "90.52", "4.10", "0.0377", 3),
    @(12, "010661", "安信稳健聚申一年持有期混合C",             "2.28",   "39.19", "1.48", "0.0337", 8),
    @(13, "012257", "安信丰穗一年持有混合C",                   "2.43",   "24.64", "1.08", "0.0262", 9),
    @(14, "011652", "招商港股通核心精选股票C",                 "0.88",   "84.28", "2.80", "0.0246", 7),
    @(15, "010754", "招商沪港深科技创新主题精选灵活配置混合C", "0.25",   "90.52", "4.10", "0.0102", 3),
    @(16, "004532", "民生加银中证港股通高股息精选指数A",       "0.13",   "92.87", "3.66", "0.0048", 9),
    @(17, "004533", "民生加银中证港股通高股息精选指数C",       "0.08",   "92.87", "3.66", "0.0029", 9)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    # Columns B, D, E, F, G look numeric but must stay text (leading zeros /
    # trailing zeros matter), so force them in with a leading apostrophe —
    # the same thing typing them by hand in Excel would do.
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Restore the originally-active sheet/tab ("总计") so this edit doesn't
# leave an unrelated UI-state change (active tab) in the diff.
$summary.Activate()

Write-Host "Edit complete"
